# Refresh the cryptocurrency price/volume snapshot on the active sheet.
# Each data row layout is: A=index, B=Coin, C=Link, D=Price, E=Volume(1h).
# Prices that are stored as plain text in the sheet (e.g. "64.069.42",
# "0.598") are written back with a leading apostrophe so Excel keeps them
# as text instead of re-interpreting them as numbers/dates.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
    $ws.Cells.Item(2, 4).Value = '64.069.42'   # Price
    $ws.Cells.Item(2, 5).Value = '  -3.41%  '   # Volume(1h)

# Row 3
    $ws.Cells.Item(3, 4).Value = '3.151.74'   # Price
    $ws.Cells.Item(3, 5).Value = '  -4.98%  '   # Volume(1h)

# Row 4
    $ws.Cells.Item(4, 5).Value = '  +0.09%  '   # Volume(1h)

# Row 5
    $ws.Cells.Item(5, 4).Value = '''567.75'   # Price (kept as text)
    $ws.Cells.Item(5, 5).Value = '  -3.03%  '   # Volume(1h)

# Row 6
    $ws.Cells.Item(6, 4).Value = '''166.29'   # Price (kept as text)
    $ws.Cells.Item(6, 5).Value = '  -7.77%  '   # Volume(1h)

# Row 7
    $ws.Cells.Item(7, 4).Value = '''0.598'   # Price (kept as text)
    $ws.Cells.Item(7, 5).Value = '  -8.37%  '   # Volume(1h)

# Row 8
    $ws.Cells.Item(8, 5).Value = '  +0.05%  '   # Volume(1h)

# Row 9
    $ws.Cells.Item(9, 4).Value = '3.154.74'   # Price
    $ws.Cells.Item(9, 5).Value = '  -4.85%  '   # Volume(1h)

# Row 10
    $ws.Cells.Item(10, 5).Value = '  -7.07%  '   # Volume(1h)

# Row 11
    $ws.Cells.Item(11, 4).Value = '''6.77'   # Price (kept as text)
    $ws.Cells.Item(11, 5).Value = '  -0.68%  '   # Volume(1h)

# Row 12
    $ws.Cells.Item(12, 4).Value = '''0.384'   # Price (kept as text)
    $ws.Cells.Item(12, 5).Value = '  -4.33%  '   # Volume(1h)

# Row 13
    $ws.Cells.Item(13, 4).Value = '3.703.23'   # Price
    $ws.Cells.Item(13, 5).Value = '  -5.06%  '   # Volume(1h)

# Row 14
    $ws.Cells.Item(14, 5).Value = '  -1.40%  '   # Volume(1h)

# Row 15
    $ws.Cells.Item(15, 4).Value = '64.176.49'   # Price
    $ws.Cells.Item(15, 5).Value = '  -3.34%  '   # Volume(1h)

# Row 16
    $ws.Cells.Item(16, 4).Value = '''25.04'   # Price (kept as text)
    $ws.Cells.Item(16, 5).Value = '  -5.49%  '   # Volume(1h)

# Row 17
    $ws.Cells.Item(17, 4).Value = '3.159.67'   # Price
    $ws.Cells.Item(17, 5).Value = '  -3.84%  '   # Volume(1h)

# Row 18
    $ws.Cells.Item(18, 4).Value = '''0.0000155'   # Price (kept as text)
    $ws.Cells.Item(18, 5).Value = '  -5.39%  '   # Volume(1h)

# Row 19
    $ws.Cells.Item(19, 4).Value = '''414.07'   # Price (kept as text)

# Row 20
    $ws.Cells.Item(20, 4).Value = '''12.68'   # Price (kept as text)
    $ws.Cells.Item(20, 5).Value = '  -3.41%  '   # Volume(1h)

# Row 21
    $ws.Cells.Item(21, 4).Value = '''5.23'   # Price (kept as text)
    $ws.Cells.Item(21, 5).Value = '  -4.49%  '   # Volume(1h)

# Row 22
    $ws.Cells.Item(22, 4).Value = '''7.08'   # Price (kept as text)
    $ws.Cells.Item(22, 5).Value = '  -3.57%  '   # Volume(1h)

# Row 23
    $ws.Cells.Item(23, 4).Value = '''0.998'   # Price (kept as text)
    $ws.Cells.Item(23, 5).Value = '  -0.35%  '   # Volume(1h)

# Row 24
    $ws.Cells.Item(24, 5).Value = '  -0.06%  '   # Volume(1h)

# Row 25
    $ws.Cells.Item(25, 4).Value = '''69.41'   # Price (kept as text)
    $ws.Cells.Item(25, 5).Value = '  -3.18%  '   # Volume(1h)

# Row 26
    $ws.Cells.Item(26, 4).Value = '''0.204'   # Price (kept as text)
    $ws.Cells.Item(26, 5).Value = '  -0.45%  '   # Volume(1h)

# Row 27
    $ws.Cells.Item(27, 4).Value = '''0.494'   # Price (kept as text)
    $ws.Cells.Item(27, 5).Value = '  -3.90%  '   # Volume(1h)

# Row 28
    $ws.Cells.Item(28, 4).Value = '''0.0000101'   # Price (kept as text)
    $ws.Cells.Item(28, 5).Value = '  -12.21%  '   # Volume(1h)

# Row 29
    $ws.Cells.Item(29, 4).Value = '''8.71'   # Price (kept as text)
    $ws.Cells.Item(29, 5).Value = '  -4.24%  '   # Volume(1h)

# Row 30
    $ws.Cells.Item(30, 4).Value = '''0.999'   # Price (kept as text)
    $ws.Cells.Item(30, 5).Value = '  +0.05%  '   # Volume(1h)

# Row 31
    $ws.Cells.Item(31, 5).Value = '  -5.24%  '   # Volume(1h)

# Row 32
    $ws.Cells.Item(32, 2).Value = 'USDe'   # Coin
    $ws.Cells.Item(32, 3).Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'   # Link
    $ws.Cells.Item(32, 4).Value = '''0.999'   # Price (kept as text)
    $ws.Cells.Item(32, 5).Value = '  -0.12%  '   # Volume(1h)

# Row 33
    $ws.Cells.Item(33, 2).Value = 'EthereumClassic'   # Coin
    $ws.Cells.Item(33, 3).Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'   # Link
    $ws.Cells.Item(33, 4).Value = '''21.55'   # Price (kept as text)
    $ws.Cells.Item(33, 5).Value = '  -3.61%  '   # Volume(1h)

# Row 34
    $ws.Cells.Item(34, 4).Value = '''4.97'   # Price (kept as text)
    $ws.Cells.Item(34, 5).Value = '  -4.02%  '   # Volume(1h)

# Row 35
    $ws.Cells.Item(35, 4).Value = '''6.28'   # Price (kept as text)
    $ws.Cells.Item(35, 5).Value = '  -4.90%  '   # Volume(1h)

# Row 36
    $ws.Cells.Item(36, 2).Value = 'Monero'   # Coin
    $ws.Cells.Item(36, 3).Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'   # Link
    $ws.Cells.Item(36, 4).Value = '''155.44'   # Price (kept as text)
    $ws.Cells.Item(36, 5).Value = '  -2.75%  '   # Volume(1h)

# Row 37
    $ws.Cells.Item(37, 2).Value = 'Fetch.AI'   # Coin
    $ws.Cells.Item(37, 3).Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'   # Link
    $ws.Cells.Item(37, 4).Value = '''1.11'   # Price (kept as text)
    $ws.Cells.Item(37, 5).Value = '  -6.24%  '   # Volume(1h)

# Row 38
    $ws.Cells.Item(38, 4).Value = '''1.34'   # Price (kept as text)
    $ws.Cells.Item(38, 5).Value = '  -6.18%  '   # Volume(1h)

# Row 39
    $ws.Cells.Item(39, 4).Value = '2.686.10'   # Price
    $ws.Cells.Item(39, 5).Value = '  -6.24%  '   # Volume(1h)

# Row 40
    $ws.Cells.Item(40, 4).Value = '''1.67'   # Price (kept as text)
    $ws.Cells.Item(40, 5).Value = '  -6.96%  '   # Volume(1h)

# Row 41
    $ws.Cells.Item(41, 4).Value = '''4.15'   # Price (kept as text)
    $ws.Cells.Item(41, 5).Value = '  -3.88%  '   # Volume(1h)

# Row 42
    $ws.Cells.Item(42, 4).Value = '''23.78'   # Price (kept as text)
    $ws.Cells.Item(42, 5).Value = '  -9.70%  '   # Volume(1h)

# Row 43
    $ws.Cells.Item(43, 2).Value = 'Mantle'   # Coin
    $ws.Cells.Item(43, 3).Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'   # Link
    $ws.Cells.Item(43, 4).Value = '''0.717'   # Price (kept as text)
    $ws.Cells.Item(43, 5).Value = '  -5.49%  '   # Volume(1h)

# Row 44
    $ws.Cells.Item(44, 2).Value = 'OKB'   # Coin
    $ws.Cells.Item(44, 3).Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'   # Link
    $ws.Cells.Item(44, 4).Value = '''39.01'   # Price (kept as text)
    $ws.Cells.Item(44, 5).Value = '  -1.77%  '   # Volume(1h)

# Row 45
    $ws.Cells.Item(45, 4).Value = '''0.0611'   # Price (kept as text)
    $ws.Cells.Item(45, 5).Value = '  -7.33%  '   # Volume(1h)

# Row 46
    $ws.Cells.Item(46, 4).Value = '''5.41'   # Price (kept as text)
    $ws.Cells.Item(46, 5).Value = '  -8.32%  '   # Volume(1h)

# Row 47
    $ws.Cells.Item(47, 4).Value = '''0.0259'   # Price (kept as text)
    $ws.Cells.Item(47, 5).Value = '  -4.70%  '   # Volume(1h)

# Row 48
    $ws.Cells.Item(48, 4).Value = '''286.89'   # Price (kept as text)
    $ws.Cells.Item(48, 5).Value = '  -8.00%  '   # Volume(1h)

# Row 49
    $ws.Cells.Item(49, 4).Value = '''21.01'   # Price (kept as text)
    $ws.Cells.Item(49, 5).Value = '  -9.03%  '   # Volume(1h)

# Row 50
    $ws.Cells.Item(50, 5).Value = '  +0.05%  '   # Volume(1h)

# Row 51
    $ws.Cells.Item(51, 4).Value = '''0.0983'   # Price (kept as text)
    $ws.Cells.Item(51, 5).Value = '  -5.98%  '   # Volume(1h)
